# Update the lattice-multiplication exercise table to the new set of
# problems / partial-product digits, cell by cell. Each cell's text is
# laid out as 5 lines separated by Word line-breaks (vertical tab, char 11,
# which Word's Range.Text setter turns into <w:br/> elements):
#   <A> x <B>
#     <tens of B>    <ones of B>
#   ----
#   <tens of A>|    |
#   <ones of A>|    |

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

$t.Cell(1,1).Range.Text = "82 x 20" + $nl + "  2    0" + $nl + "  ----" + $nl + "8|    |" + $nl + "2|    |"
$t.Cell(1,2).Range.Text = "24 x 29" + $nl + "  2    9" + $nl + "  ----" + $nl + "2|    |" + $nl + "4|    |"
$t.Cell(1,3).Range.Text = "81 x 19" + $nl + "  1    9" + $nl + "  ----" + $nl + "8|    |" + $nl + "1|    |"

$t.Cell(2,1).Range.Text = "53 x 60" + $nl + "  6    0" + $nl + "  ----" + $nl + "5|    |" + $nl + "3|    |"
$t.Cell(2,2).Range.Text = "55 x 37" + $nl + "  3    7" + $nl + "  ----" + $nl + "5|    |" + $nl + "5|    |"
$t.Cell(2,3).Range.Text = "20 x 93" + $nl + "  9    3" + $nl + "  ----" + $nl + "2|    |" + $nl + "0|    |"

$t.Cell(3,1).Range.Text = "88 x 48" + $nl + "  4    8" + $nl + "  ----" + $nl + "8|    |" + $nl + "8|    |"
$t.Cell(3,2).Range.Text = "29 x 86" + $nl + "  8    6" + $nl + "  ----" + $nl + "2|    |" + $nl + "9|    |"
$t.Cell(3,3).Range.Text = "33 x 29" + $nl + "  2    9" + $nl + "  ----" + $nl + "3|    |" + $nl + "3|    |"

$t.Cell(4,1).Range.Text = "81 x 10" + $nl + "  1    0" + $nl + "  ----" + $nl + "8|    |" + $nl + "1|    |"
$t.Cell(4,2).Range.Text = "22 x 88" + $nl + "  8    8" + $nl + "  ----" + $nl + "2|    |" + $nl + "2|    |"
$t.Cell(4,3).Range.Text = "73 x 68" + $nl + "  6    8" + $nl + "  ----" + $nl + "7|    |" + $nl + "3|    |"

$t.Cell(5,1).Range.Text = "91 x 62" + $nl + "  6    2" + $nl + "  ----" + $nl + "9|    |" + $nl + "1|    |"
$t.Cell(5,2).Range.Text = "85 x 47" + $nl + "  4    7" + $nl + "  ----" + $nl + "8|    |" + $nl + "5|    |"
$t.Cell(5,3).Range.Text = "40 x 80" + $nl + "  8    0" + $nl + "  ----" + $nl + "4|    |" + $nl + "0|    |"
